$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 41.43823833333334
$ws.Range("H2").Value = 124.314715
$ws.Range("I2").Value = 0.981992391336623
$ws.Range("J2").Value = 0.9819923913366232
$ws.Range("M2").Value = 7.413580666666667
$ws.Range("N2").Value = 22.240742
$ws.Range("O2").Value = 0.05108888817597561
$ws.Range("P2").Value = 0.05108888817597561
$ws.Range("Q2").Value = 307.2057225687256
$ws.Range("R2").Value = 2764.85150311853
$ws.Range("S2").Value = 0.05016889947065561
$ws.Range("T2").Value = 0.05016889947065562

# Row 3
$ws.Range("G3").Value = 41.43823833333334
$ws.Range("H3").Value = 124.314715
$ws.Range("I3").Value = 0.981992391336623
$ws.Range("J3").Value = 0.9819923913366232
$ws.Range("O3").Value = 0.0112127179963522
$ws.Range("P3").Value = 0.0112127179963522
$ws.Range("Q3").Value = 67.42388133724445
$ws.Range("R3").Value = 606.8149320352001
$ws.Range("S3").Value = 0.01101080375862109
$ws.Range("T3").Value = 0.01101080375862109

# Row 4
$ws.Range("G4").Value = 41.43823833333334
$ws.Range("H4").Value = 124.314715
$ws.Range("I4").Value = 0.981992391336623
$ws.Range("J4").Value = 0.9819923913366232
$ws.Range("O4").Value = 0.9376983938276722
$ws.Range("P4").Value = 0.9376983938276722
$ws.Range("Q4").Value = 5638.531643811066
$ws.Range("R4").Value = 50746.78479429959
$ws.Range("S4").Value = 0.9208126881073463
$ws.Range("T4").Value = 0.9208126881073465

# Row 5
$ws.Range("I5").Value = 0.006845967574057415
$ws.Range("J5").Value = 0.006845967574057417
$ws.Range("M5").Value = 7.413580666666667
$ws.Range("N5").Value = 22.240742
$ws.Range("O5").Value = 0.05108888817597561
$ws.Range("P5").Value = 0.05108888817597561
$ws.Range("Q5").Value = 2.141687078051334
$ws.Range("R5").Value = 19.275183702462
$ws.Range("S5").Value = 0.0003497528718473743
$ws.Range("T5").Value = 0.0003497528718473744

# Row 6
$ws.Range("I6").Value = 0.006845967574057415
$ws.Range("J6").Value = 0.006845967574057417
$ws.Range("O6").Value = 0.0112127179963522
$ws.Range("P6").Value = 0.0112127179963522
$ws.Range("S6").Value = 0.00007676190382007719
$ws.Range("T6").Value = 0.00007676190382007722

# Row 7
$ws.Range("I7").Value = 0.006845967574057415
$ws.Range("J7").Value = 0.006845967574057417
$ws.Range("O7").Value = 0.9376983938276722
$ws.Range("P7").Value = 0.9376983938276722
$ws.Range("S7").Value = 0.006419452798389963
$ws.Range("T7").Value = 0.006419452798389965

# Row 8
$ws.Range("I8").Value = 0.01116164108931947
$ws.Range("J8").Value = 0.01116164108931947
$ws.Range("M8").Value = 7.413580666666667
$ws.Range("N8").Value = 22.240742
$ws.Range("O8").Value = 0.05108888817597561
$ws.Range("P8").Value = 0.05108888817597561
$ws.Range("Q8").Value = 3.491798965193555
$ws.Range("R8").Value = 31.426190686742
$ws.Range("S8").Value = 0.000570235833472617
$ws.Range("T8").Value = 0.0005702358334726172

# Row 9
$ws.Range("I9").Value = 0.01116164108931947
$ws.Range("J9").Value = 0.01116164108931947
$ws.Range("O9").Value = 0.0112127179963522
$ws.Range("P9").Value = 0.0112127179963522
$ws.Range("S9").Value = 0.0001251523339110366
$ws.Range("T9").Value = 0.0001251523339110366

# Row 10
$ws.Range("I10").Value = 0.01116164108931947
$ws.Range("J10").Value = 0.01116164108931947
$ws.Range("O10").Value = 0.9376983938276722
$ws.Range("P10").Value = 0.9376983938276722
$ws.Range("S10").Value = 0.01046625292193581
$ws.Range("T10").Value = 0.01046625292193582

